# Auto-generated edit script applying the cryptos.xlsx "Updated cryptos list" diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") always holds text in this workbook, even when the text looks
# like a plain decimal number (e.g. "2.31"). Force every Price cell we touch to the
# Text number format first so Excel does not silently convert it to a numeric value.
$priceCells = @('D2', 'D3', 'D5', 'D6', 'D7', 'D10', 'D11', 'D15', 'D16', 'D17', 'D18', 'D19', 'D22', 'D23', 'D24', 'D25', 'D26', 'D28', 'D29', 'D30', 'D31', 'D32', 'D33', 'D36', 'D37', 'D38', 'D39', 'D41', 'D43', 'D45', 'D49', 'D51')
foreach ($cell in $priceCells) {
    $ws.Range($cell).NumberFormat = '@'
}

$ws.Range('D2').Value = '48.283.51'
$ws.Range('E2').Value = '  +1.86%  '
$ws.Range('D3').Value = '2.522.09'
$ws.Range('E3').Value = '  +0.69%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = '323.08'
$ws.Range('E5').Value = '  -0.51%  '
$ws.Range('D6').Value = '109.29'
$ws.Range('E6').Value = '  -0.74%  '
$ws.Range('D7').Value = '0.528'
$ws.Range('E7').Value = '  +0.49%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E9').Value = '  +4.23%  '
$ws.Range('D10').Value = '40.35'
$ws.Range('E10').Value = '  +2.34%  '
$ws.Range('D11').Value = '20.02'
$ws.Range('E11').Value = '  +7.92%  '
$ws.Range('E12').Value = '  +0.40%  '
$ws.Range('E13').Value = '  +0.73%  '
$ws.Range('E14').Value = '  +0.40%  '
$ws.Range('D15').Value = '2.914.09'
$ws.Range('E15').Value = '  +0.57%  '
$ws.Range('D16').Value = '2.508.46'
$ws.Range('E16').Value = '  +0.00%  '
$ws.Range('D17').Value = '0.860'
$ws.Range('E17').Value = '  -0.08%  '
$ws.Range('D18').Value = '48.162.57'
$ws.Range('E18').Value = '  +1.73%  '
$ws.Range('D19').Value = '13.31'
$ws.Range('E19').Value = '  +3.20%  '
$ws.Range('E20').Value = '  -0.41%  '
$ws.Range('E21').Value = '  +0.20%  '
$ws.Range('D22').Value = '2.73'
$ws.Range('E22').Value = '  +0.68%  '
$ws.Range('D23').Value = '72.52'
$ws.Range('E23').Value = '  +2.73%  '
$ws.Range('D24').Value = '268.38'
$ws.Range('E24').Value = '  +7.84%  '
$ws.Range('D25').Value = '2.58'
$ws.Range('E25').Value = '  -1.63%  '
$ws.Range('D26').Value = '26.19'
$ws.Range('E26').Value = '  +0.28%  '
$ws.Range('E27').Value = '  +0.12%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').Value = '2.31'
$ws.Range('E28').Value = '  +2.22%  '
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').Value = '10.20'
$ws.Range('E29').Value = '  +1.25%  '
$ws.Range('B30').Value = 'Kaspa'
$ws.Range('C30').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D30').Value = '0.145'
$ws.Range('E30').Value = '  +5.24%  '
$ws.Range('D31').Value = '35.14'
$ws.Range('E31').Value = '  -1.04%  '
$ws.Range('D32').Value = '49.82'
$ws.Range('E32').Value = '  -0.24%  '
$ws.Range('D33').Value = '20.02'
$ws.Range('E33').Value = '  -0.06%  '
$ws.Range('E34').Value = '  -0.37%  '
$ws.Range('E35').Value = '  -0.07%  '
$ws.Range('D36').Value = '0.0792'
$ws.Range('E36').Value = '  -0.78%  '
$ws.Range('D37').Value = '1.99'
$ws.Range('E37').Value = '  -0.53%  '
$ws.Range('D38').Value = '4.74'
$ws.Range('E38').Value = '  +0.63%  '
$ws.Range('D39').Value = '3.01'
$ws.Range('E39').Value = '  +0.29%  '
$ws.Range('E40').Value = '  +0.11%  '
$ws.Range('D41').Value = '22.29'
$ws.Range('E41').Value = '  +4.53%  '
$ws.Range('E42').Value = '  -1.34%  '
$ws.Range('D43').Value = '118.32'
$ws.Range('E43').Value = '  -3.05%  '
$ws.Range('E44').Value = '  +0.03%  '
$ws.Range('D45').Value = '2.002.75'
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('E46').Value = '  -0.27%  '
$ws.Range('E48').Value = '  -2.08%  '
$ws.Range('D49').Value = '9.11'
$ws.Range('E49').Value = '  +0.45%  '
$ws.Range('E50').Value = '  +0.51%  '
$ws.Range('D51').Value = '80.68'
$ws.Range('E51').Value = '  +3.03%  '
